$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.797.50"
$ws.Range("E2").Value = "  -3.70%  "

$ws.Range("D3").Value = "2.906.96"
$ws.Range("E3").Value = "  -4.21%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "585.33"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.29%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.81"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -6.14%  "

$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("E8").Value = "  -2.59%  "

$ws.Range("D9").Value = "2.908.08"
$ws.Range("E9").Value = "  -4.07%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.69"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.51%  "

$ws.Range("E11").Value = "  -4.73%  "

$ws.Range("E12").Value = "  -3.95%  "

$ws.Range("E13").Value = "  -3.40%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.40"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -6.60%  "

$ws.Range("E15").Value = "  +1.58%  "

$ws.Range("D16").Value = "3.387.29"
$ws.Range("E16").Value = "  -4.26%  "

$ws.Range("D17").Value = "60.744.71"
$ws.Range("E17").Value = "  -3.64%  "

$ws.Range("E18").Value = "  -5.45%  "

$ws.Range("D19").Value = "2.904.93"
$ws.Range("E19").Value = "  -4.24%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "428.81"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -5.41%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.56"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.99%  "

$ws.Range("E22").Value = "  -2.40%  "

$ws.Range("E23").Value = "  -5.34%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "80.63"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.01%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "10.95"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.49%  "

$ws.Range("E26").Value = "  -2.89%  "

$ws.Range("E27").Value = "  -4.31%  "

$ws.Range("E28").Value = "  +0.00%  "

$ws.Range("E29").Value = "  -2.88%  "

$ws.Range("E30").Value = "  +0.07%  "

$ws.Range("E31").Value = "  -3.59%  "

$ws.Range("E32").Value = "  -3.38%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "26.49"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.94%  "

$ws.Range("E34").Value = "  -3.40%  "

$ws.Range("E35").Value = "  +1.45%  "

$ws.Range("E36").Value = "  -3.32%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.61"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.36%  "

$ws.Range("E38").Value = "  -4.43%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "49.53"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.95%  "

$ws.Range("E40").Value = "  -4.43%  "

$ws.Range("E41").Value = "  -4.74%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.60"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.78%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.296"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.91%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "41.48"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.63%  "

$ws.Range("E45").Value = "  -2.63%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "378.02"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.43%  "

$ws.Range("D47").Value = "2.700.94"
$ws.Range("E47").Value = "  -0.78%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "132.33"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.80%  "

$ws.Range("E49").Value = "  -0.02%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "24.24"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.60%  "

$ws.Range("E51").Value = "  -2.63%  "
